$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.9
$ws.Range("J2").Value = 3.9
$ws.Range("Q2").Value = 1.58
$ws.Range("S2").Value = 2.44
$ws.Range("V2").Value = 1.63
$ws.Range("W2").Value = 1.52
$ws.Range("Z2").Value = 22
$ws.Range("AH2").Value = 14
$ws.Range("AJ2").Value = 46
$ws.Range("H3").Value = 5.2
$ws.Range("I3").Value = 6.2
$ws.Range("K3").Value = 6.2
$ws.Range("P3").Value = 3.4
$ws.Range("R3").Value = 1.99
$ws.Range("S3").Value = 1.82
$ws.Range("U3").Value = 2.72
$ws.Range("V3").Value = 1.2
$ws.Range("W3").Value = 2.66
$ws.Range("Y3").Value = 48
$ws.Range("AB3").Value = 19.5
$ws.Range("AC3").Value = 15
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 15.5
$ws.Range("F4").Value = 1.52
$ws.Range("K4").Value = 4.7
$ws.Range("O4").Value = 1.33
$ws.Range("U4").Value = 1.8
$ws.Range("F5").Value = 2.76
$ws.Range("J5").Value = 2.98
$ws.Range("K5").Value = 3.35
$ws.Range("R5").Value = 1.22
$ws.Range("U5").Value = 1.86
$ws.Range("V5").Value = 1.48
$ws.Range("G6").Value = 3.8
$ws.Range("W6").Value = 1.36
$ws.Range("G10").Value = 2.48
$ws.Range("H10").Value = 3.35
$ws.Range("F11").Value = 1.26
$ws.Range("G11").Value = 1.33
$ws.Range("J11").Value = 5.6
$ws.Range("L11").Value = 1.17
$ws.Range("N11").Value = 7.6
$ws.Range("O11").Value = 1.09
$ws.Range("P11").Value = 3.3
$ws.Range("Q11").Value = 1.28
$ws.Range("R11").Value = 2.08
$ws.Range("S11").Value = 1.7
$ws.Range("T11").Value = 1.64
$ws.Range("U11").Value = 2.2
$ws.Range("V11").Value = 1.07
$ws.Range("W11").Value = 3.65
$ws.Range("AN11").Value = 3.85
$ws.Range("G12").Value = 2.94
$ws.Range("I12").Value = 2.6
$ws.Range("W12").Value = 1.51
$ws.Range("AA12").Value = 1000
$ws.Range("AE12").Value = 1000
$ws.Range("AN12").Value = 17
$ws.Range("G15").Value = 4.6
$ws.Range("R15").Value = 1.23
$ws.Range("V15").Value = 1.75
$ws.Range("F16").Value = 3.5
$ws.Range("H16").Value = 2.58
$ws.Range("I16").Value = 2.68
$ws.Range("K16").Value = 2.98
$ws.Range("AA16").Value = 46
$ws.Range("AE16").Value = 55
$ws.Range("F18").Value = 4.2
$ws.Range("W18").Value = 1.21
$ws.Range("G19").Value = 1.95
$ws.Range("V19").Value = 1.23
$ws.Range("AB19").Value = 8.4
